$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns (Wins, Losses, Ties) after the existing
# "Unnamed: 28" column (AC), matching the header style used by the rest of
# row 1 (bold, centered, bordered -> same style as AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Populate the season record (Wins / Losses / Ties) for every player row.
$wins = 83
$losses = 78
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
